$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.086.56'
$ws.Range('E2').Value = '  +0.00%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.651.32'
$ws.Range('E3').Value = '  +0.02%  '
$ws.Range('E4').Value = '  -0.34%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.17'
$ws.Range('E5').Value = '  +0.34%  '
$ws.Range('E6').Value = '  -0.17%  '
$ws.Range('E7').Value = '  -0.33%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2641'
$ws.Range('E8').Value = '  +0.95%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06334'
$ws.Range('E9').Value = '  +0.69%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.39'
$ws.Range('E10').Value = '  -0.54%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07692'
$ws.Range('E11').Value = '  -1.27%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.591'
$ws.Range('E12').Value = '  +2.53%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.649.57'
$ws.Range('E13').Value = '  +3.87%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.878.52'
$ws.Range('E15').Value = '  +0.84%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₅8146'
$ws.Range('E16').Value = '  +1.75%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.37'
$ws.Range('E17').Value = '  +0.90%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '26.099.40'
$ws.Range('E18').Value = '  +0.08%  '
$ws.Range('E19').Value = '  -0.36%  '
$ws.Range('E20').Value = '  -0.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.47'
$ws.Range('E21').Value = '  +4.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '191.48'
$ws.Range('E22').Value = '  -1.33%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.924'
$ws.Range('E23').Value = '  -0.34%  '
$ws.Range('E24').Value = '  -0.32%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.17'
$ws.Range('E25').Value = '  -1.63%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1191'
$ws.Range('E26').Value = '  -0.98%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.221'
$ws.Range('E28').Value = '  +0.29%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.503'
$ws.Range('E29').Value = '  +1.79%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05491'
$ws.Range('E30').Value = '  -1.82%  '
$ws.Range('E31').Value = '  +0.11%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.439'
$ws.Range('E32').Value = '  -1.07%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.357'
$ws.Range('E33').Value = '  +0.39%  '
$ws.Range('E34').Value = '  -2.32%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.420'
$ws.Range('E35').Value = '  +0.44%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9479'
$ws.Range('E36').Value = '  -0.22%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.783'
$ws.Range('E37').Value = '  -0.55%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5634'
$ws.Range('E38').Value = '  -0.17%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01579'
$ws.Range('E39').Value = '  -0.23%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.848'
$ws.Range('E40').Value = '  -1.74%  '
$ws.Range('E41').Value = '  -0.26%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.027.85'
$ws.Range('E42').Value = '  -2.74%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8292'
$ws.Range('E43').Value = '  -1.40%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.19'
$ws.Range('E44').Value = '  -1.35%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.792.21'
$ws.Range('E45').Value = '  +0.17%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '57.49'
$ws.Range('E46').Value = '  +0.54%  '
$ws.Range('E47').Value = '  +3.77%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.9983'
$ws.Range('E48').Value = '  -0.97%  '
$ws.Range('E49').Value = '  -0.04%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.975'
$ws.Range('E50').Value = '  +0.40%  '
$ws.Range('E51').Value = '  -2.64%  '
